$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '28.129.94'
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.825.95'
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  +1.36%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '311.86'
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  -0.80%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  -0.34%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.5116'
$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  -2.72%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3944'
$ws.Range("E8").NumberFormat = '@'
$ws.Range("E8").Value = '  +3.18%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.1007'
$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  +25.86%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '40.99'

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '6.485'
$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  +2.77%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  -0.35%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '20.69'
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("B15").NumberFormat = '@'
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").NumberFormat = '@'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '1.832.40'

$ws.Range("B16").NumberFormat = '@'
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").NumberFormat = '@'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '7.399'
$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  +1.11%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '94.85'
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  +3.02%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.00001138'
$ws.Range("E18").NumberFormat = '@'
$ws.Range("E18").Value = '  +4.13%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.06602'
$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.060'
$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  +1.48%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '28.205.78'
$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '11.20'
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  +0.52%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.241'
$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  -1.15%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.476'
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  +5.14%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '158.93'
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  -1.23%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '20.83'
$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  +1.85%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '2.034.10'
$ws.Range("E29").NumberFormat = '@'
$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '128.70'
$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  +4.63%  '

$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.071'
$ws.Range("E32").NumberFormat = '@'
$ws.Range("E32").Value = '  +1.34%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '5.644'
$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '3.632'
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  -1.49%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.06916'
$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  -4.79%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '9.144'
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  +6.25%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.02349'
$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.2173'
$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  +1.20%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '11.63'
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  -5.63%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '5.030'
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  -1.29%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.6271'
$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  +1.26%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.0000'
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '1.162'
$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  -0.38%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '13.37'
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  +0.65%  '

$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '3.717'
$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  -1.44%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.289'
$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  -5.96%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '125.91'
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  -1.01%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.990'
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  +3.43%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.190'
$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  -2.94%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.06783'
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -0.28%  '
